# Add new legal-* namespace rows (HK, JP, KR, MO, MY, PH, SG, TH, TW) to the
# "Namespaces" worksheet, inserted right after the existing "legal-us" block
# (i.e. directly above the "sector-education" row), matching the upstream
# commit "update namespaces for legal extensions".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Namespaces")

$firstNewRow = 56
$newCount = 9

# Push the existing rows (sector-education ... p7012) down by $newCount rows.
$ws.Rows($firstNewRow.ToString() + ":" + ($firstNewRow + $newCount - 1).ToString()).Insert()

# New legal namespace entries: prefix, namespace (fragment incl. trailing #),
# comment, and the hyperlink target (namespace URL without the trailing #).
$entries = @(
    @("legal-hk", "https://w3id.org/dpv/legal/hk#", "Laws and Authorities for HK", "https://w3id.org/dpv/legal/hk"),
    @("legal-jp", "https://w3id.org/dpv/legal/jp#", "Laws and Authorities for JP", "https://w3id.org/dpv/legal/jp"),
    @("legal-kr", "https://w3id.org/dpv/legal/kr#", "Laws and Authorities for KR", "https://w3id.org/dpv/legal/kr"),
    @("legal-mo", "https://w3id.org/dpv/legal/mo#", "Laws and Authorities for MO", "https://w3id.org/dpv/legal/mo"),
    @("legal-my", "https://w3id.org/dpv/legal/my#", "Laws and Authorities for MY", "https://w3id.org/dpv/legal/my"),
    @("legal-ph", "https://w3id.org/dpv/legal/ph#", "Laws and Authorities for PH", "https://w3id.org/dpv/legal/ph"),
    @("legal-sg", "https://w3id.org/dpv/legal/sg#", "Laws and Authorities for SG", "https://w3id.org/dpv/legal/sg"),
    @("legal-th", "https://w3id.org/dpv/legal/th#", "Laws and Authorities for TH", "https://w3id.org/dpv/legal/th"),
    @("legal-tw", "https://w3id.org/dpv/legal/tw#", "Laws and Authorities for TW", "https://w3id.org/dpv/legal/tw")
)

# All these rows share the same "approved" review date: 2025-03-20.
$reviewDate = 45736

for ($i = 0; $i -lt $entries.Count; $i++) {
    $r = $firstNewRow + $i
    $prefix = $entries[$i][0]
    $namespaceUri = $entries[$i][1]
    $comment = $entries[$i][2]
    $hyperlinkTarget = $entries[$i][3]

    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellC = $ws.Cells.Item($r, 3)
    $cellD = $ws.Cells.Item($r, 4)
    $cellE = $ws.Cells.Item($r, 5)

    $cellA.Value = $prefix
    $cellB.Value = $namespaceUri
    $cellC.Value = $comment
    $cellD.Value = $reviewDate
    $cellD.NumberFormat = "yyyy-mm-dd"
    $cellE.Value = "approved"

    $ws.Hyperlinks.Add($cellB, $hyperlinkTarget) | Out-Null

    # Plain text styling (matches the rest of the "approved" rows: Arial 10, black).
    $cellA.Font.Name = "Arial"
    $cellA.Font.Size = 10
    $cellA.Font.Color = 0
    $cellA.Font.Underline = 0

    $cellC.Font.Name = "Arial"
    $cellC.Font.Size = 10
    $cellC.Font.Color = 0
    $cellC.Font.Underline = 0

    $cellE.Font.Name = "Arial"
    $cellE.Font.Size = 10
    $cellE.Font.Color = 0
    $cellE.Font.Underline = 0

    # Hyperlink-style text for the namespace URI column.
    $cellB.Font.Name = "Arial"
    $cellB.Font.Size = 10
    $cellB.Font.Color = 16711680
    $cellB.Font.Underline = 2
}
